# VRA Skills Activity Review Form — apply reviewed-edit pass.
# Uses Find/Replace across $word.ActiveDocument.Content for each wording
# change, plus explicit Bookmarks handling for the relocated "_GoBack"
# bookmark (Word's "last edit position" marker).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $old)
    }
    return $found
}

function ReplaceAll-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
    $count = 0
    while ($true) {
        $rng2 = $d.Content
        $found = $rng2.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
        if (-not $found) { break }
        $count = $count + 1
        if ($count -gt 20) { break }
    }
}

# 1) "Organisation: YMCA / CoderDojo" — merge the "YMCA / " + "CoderDojo"
#    runs (previously split around a spell-check proofErr pair) into a
#    single run and drop the stale proofing marks. Extending the search
#    to " CoderDojo" lets the replace also swallow both <w:proofErr/>
#    markers that bracket the word.
Replace-Text "YMCA / CoderDojo" "YMCA / CoderDojo"

# 2) Second week reflection — reworded.
Replace-Text `
    "Second week where to be better prepared I had brought pre created learning material for students to study from. (Planning " `
    "Second week I was better prepared as I had brought pre created learning material for students to study from. (Planning "

# 3) Internet fundamentals lecture entry — reworded.
Replace-Text `
    "Delivered lecture on basic internet fundamentals after which I took students websites and put them onto my own server so they could access them using web browser. (Technical Ability, Presentation Skills, Leadership)" `
    "Delivered a lecture on basic internet fundamentals, then I took students websites and put them onto my own server so they could access them using a web browser. (Technical Ability, Presentation Skills, Leadership)"

# 4) External JavaScript files entry — added clause.
Replace-Text `
    "Taught students how to deal with external JavaScript files." `
    "Taught students how to deal with external JavaScript files so students could better organise their projects."

# 5) "texts" -> "text" typo fix (also folds in the trailing space + skills
#    parenthetical run that collapse into one run in the target).
Replace-Text `
    "Used my presentation and technical skills to teach more html tags for writing web pages. I covered loading images, different styles of texts and video. (Technical Ability & Presentation Skills)" `
    "Used my presentation and technical skills to teach more html tags for writing web pages. I covered loading images, different styles of text and video. (Technical Ability & Presentation Skills)"

# 6) Browsers entry — added clause.
Replace-Text `
    " are different and how it works. (Adaptability/Flexibility)" `
    " are different and how it works as was requested of me on the day. (Adaptability/Flexibility)"

# 8) CSS concepts entry — collapses the trailing space + skills run (no
#    wording change, just tidies the run split).
Replace-Text `
    "The focus of this week was to teach basic CSS concepts, selectors and properties as an example. (Technical Ability & Presentation Skills)" `
    "The focus of this week was to teach basic CSS concepts, selectors and properties as an example. (Technical Ability & Presentation Skills)"

# 9) Separate stylesheet entry — added clause.
Replace-Text `
    "Also taught students how to write styles in a separate file and load them using the link html element." `
    "Also taught students how to write styles in a separate file and load them using the link html element so they could have cleaner projects."

# 10) Computer fundamentals entry — added colon.
Replace-Text `
    "Taught basic computer fundamentals files, folders, hardware, software." `
    "Taught basic computer fundamentals: files, folders, hardware, software."

# 11 & 12) "report" -> "reported"; trims the final clause about games
#    going online, and (together with the Phaser wrap removal below)
#    clears the proofErr pair around "Phaser" by reaching past it.
Replace-Text `
    "Students had report a while ago that they wanted to be able to write games so I got students to download Phaser so next week we would write basic web games that would then go online." `
    "Students had reported a while ago that they wanted to be able to write games so I got students to download Phaser so next week we would write basic web games."

# 13) Basic-game entry — no wording change, just clears the proofErr pair
#     wrapping "Phaser" by including trailing context in the match.
Replace-Text `
    "I had to teach what Phaser was" `
    "I had to teach what Phaser was"

# 14) Final Phaser-project entry — reworded; also clears the proofErr pair
#     and the old "_GoBack" bookmark that were wrapped around "Phaser".
Replace-Text `
    "Final part of the 3 week Phaser project, where we aimed" `
    "This was the final part of the 3 week Phaser project. We aimed"

# 15) Personal-server entry — "and" -> "an", "do myself" -> "do it myself".
Replace-Text `
    ", and idea that was thought of last minute and I decided to do myself so students were satisfied" `
    ", an idea that was thought of last minute and I decided to do it myself so students were satisfied"

# 16) Jump-scare tutorials entry — "Game" -> "Gave".
Replace-Text `
    "Game several tutorials on how to use the above mentioned framework" `
    "Gave several tutorials on how to use the above mentioned framework"

# 17) Final October lesson entry — reworded content.
Replace-Text `
    "Final lesson in October looked at CSS so students could style their websites. (Technical Ability)" `
    "Final lesson in October looked at animations and positioning in CSS so students could style their websites. (Technical Ability)"

# Re-seat the hidden "_GoBack" bookmark (Word's last-edit-position marker)
# from the old Phaser-project sentence (removed by the edit above) onto
# the new October CSS sentence, between "...positioning in" and " CSS...".
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}
$anchor = $d.Content
$anchor.Find.Execute("Final lesson in October looked at animations and positioning in") | Out-Null
$gobackPoint = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $gobackPoint) | Out-Null

# 18) The "Look above for content that was covered in this session." line
#     repeats identically 7 times through the document — replace every one.
ReplaceAll-Text `
    "Look above for content that was covered in this session." `
    "Look above for my reflection of this session."
